$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 12505226
$ws.Range("I62").Value = 41669756
$ws.Range("J62").Value = 6141.4287
$ws.Range("K62").Value = 41669756
$ws.Range("L62").Value = 6141.4287
$ws.Range("M62").Value = -41669132
$ws.Range("N62").Value = -7389.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 12505226
$ws.Range("I65").Value = 41669756
$ws.Range("J65").Value = 6141.4287
$ws.Range("K65").Value = 208348780
$ws.Range("L65").Value = 30707.1435
$ws.Range("M65").Value = -208345660
$ws.Range("N65").Value = -36947.14350000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 10681.134
$ws.Range("I74").Value = 8862.846
$ws.Range("K74").Value = 8862.846
$ws.Range("M74").Value = -7926.846

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 66672076
$ws.Range("I76").Value = 4336.8335
$ws.Range("K76").Value = 4336.8335
$ws.Range("M76").Value = -4021.8335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 10681.134
$ws.Range("I77").Value = 8862.846
$ws.Range("K77").Value = 44314.23
$ws.Range("M77").Value = -39634.23

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 66672076
$ws.Range("I79").Value = 4336.8335
$ws.Range("K79").Value = 4336.8335
$ws.Range("M79").Value = -3244.8335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3010.818
$ws.Range("I132").Value = 3237.9
$ws.Range("K132").Value = 9713.700000000001
$ws.Range("M132").Value = -7183.700000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 556995.6
$ws.Range("I135").Value = 742036.5
$ws.Range("K135").Value = 6678328.5
$ws.Range("M135").Value = -6675793.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6099.3213
$ws.Range("I122").Value = 6268.385
$ws.Range("J122").Value = 5952.8
$ws.Range("K122").Value = 18805.155
$ws.Range("L122").Value = 17858.4
$ws.Range("M122").Value = -16355.155
$ws.Range("N122").Value = -22758.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 85163.164
$ws.Range("J81").Value = 85163.164
$ws.Range("L81").Value = 85163.164
$ws.Range("N81").Value = -87285.164

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 85163.164
$ws.Range("J84").Value = 85163.164
$ws.Range("L84").Value = 255489.492
$ws.Range("N84").Value = -266097.492

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 63408.875
$ws.Range("I105").Value = 67536.13
$ws.Range("K105").Value = 67536.13
$ws.Range("M105").Value = -65789.13

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 15000
$ws.Range("J109").Value = 15000
$ws.Range("L109").Value = 15000
$ws.Range("N109").Value = -17774

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 48541.8
$ws.Range("I132").Value = 42709
$ws.Range("K132").Value = 42709
$ws.Range("M132").Value = -37649

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37349.8
$ws.Range("I31").Value = 2768.6667
$ws.Range("K31").Value = 2768.6667
$ws.Range("M31").Value = -2473.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 37349.8
$ws.Range("I34").Value = 2768.6667
$ws.Range("K34").Value = 2768.6667
$ws.Range("M34").Value = -2566.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4985.328
$ws.Range("I58").Value = 4572.7085
$ws.Range("K58").Value = 4572.7085
$ws.Range("M58").Value = -4369.7085

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7950
$ws.Range("J62").Value = 7950
$ws.Range("L62").Value = 7950
$ws.Range("N62").Value = -9198

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 7950
$ws.Range("J65").Value = 7950
$ws.Range("L65").Value = 39750
$ws.Range("N65").Value = -45990

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 72000
$ws.Range("J112").Value = 72000
$ws.Range("L112").Value = 72000
$ws.Range("N112").Value = -74954

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4985.328
$ws.Range("I136").Value = 4572.7085
$ws.Range("K136").Value = 13718.1255
$ws.Range("M136").Value = -11168.1255

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 9239.25
$ws.Range("I120").Value = 9239.25
$ws.Range("K120").Value = 27717.75
$ws.Range("M120").Value = -22879.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 7004.75
$ws.Range("I133").Value = 7004.75
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 21014.25
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -15954.25
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1783.6774
$ws.Range("I140").Value = 1424.1034
$ws.Range("K140").Value = 4272.3102
$ws.Range("M140").Value = 907.6898000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9136
$ws.Range("I70").Value = 8299.556
$ws.Range("J70").Value = 11286.857
$ws.Range("K70").Value = 8299.556
$ws.Range("L70").Value = 11286.857
$ws.Range("M70").Value = -8029.556
$ws.Range("N70").Value = -11826.857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 9136
$ws.Range("I73").Value = 8299.556
$ws.Range("J73").Value = 11286.857
$ws.Range("K73").Value = 8299.556
$ws.Range("L73").Value = 11286.857
$ws.Range("M73").Value = -7363.556
$ws.Range("N73").Value = -13158.857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1431514.6
$ws.Range("I80").Value = 1003560.1
$ws.Range("J80").Value = 2501401
$ws.Range("K80").Value = 1003560.1
$ws.Range("L80").Value = 2501401
$ws.Range("M80").Value = -1002562.1
$ws.Range("N80").Value = -2503397

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1431514.6
$ws.Range("I83").Value = 1003560.1
$ws.Range("J83").Value = 2501401
$ws.Range("K83").Value = 5017800.5
$ws.Range("L83").Value = 12507005
$ws.Range("M83").Value = -5012808.5
$ws.Range("N83").Value = -12516989

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 874.6667
$ws.Range("I102").Value = 639.06665
$ws.Range("K102").Value = 639.06665
$ws.Range("M102").Value = 982.93335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 32261242
$ws.Range("I126").Value = 52634292
$ws.Range("K126").Value = 157902876
$ws.Range("M126").Value = -157900406

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 119777.664
$ws.Range("I132").Value = 13800.2
$ws.Range("K132").Value = 41400.60000000001
$ws.Range("M132").Value = -38870.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 250100000
$ws.Range("J135").Value = 250100000
$ws.Range("L135").Value = 250100000
$ws.Range("N135").Value = -250110140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 250001200
$ws.Range("I16").Value = 500000640
$ws.Range("K16").Value = 500000640
$ws.Range("M16").Value = -500000470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 45000
$ws.Range("J110").Value = 45000
$ws.Range("L110").Value = 45000
$ws.Range("N110").Value = -53180

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 405471.62
$ws.Range("I136").Value = 560936.3
$ws.Range("J136").Value = 5705.2856
$ws.Range("K136").Value = 1682808.9
$ws.Range("L136").Value = 17115.8568
$ws.Range("M136").Value = -1680258.9
$ws.Range("N136").Value = -22215.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9007
$ws.Range("J41").Value = 9249
$ws.Range("L41").Value = 9249
$ws.Range("N41").Value = -10029

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 84000
$ws.Range("J110").Value = 84000
$ws.Range("L110").Value = 84000
$ws.Range("N110").Value = -92180
